$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.760.71"
$ws.Range("E2").Value = "  +1.05%  "
$ws.Range("D3").Value = "3.637.04"
$ws.Range("E3").Value = "  +2.09%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "610.53"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.96"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.30%  "
$ws.Range("D7").Value = "3.640.03"
$ws.Range("E7").Value = "  +2.21%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("E9").Value = "  -0.23%  "
$ws.Range("E10").Value = "  +0.43%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.96"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("E12").Value = "  +1.49%  "
$ws.Range("D13").Value = "4.251.45"
$ws.Range("E13").Value = "  +2.05%  "
$ws.Range("E14").Value = "  +1.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "30.04"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.22%  "
$ws.Range("D16").Value = "3.619.56"
$ws.Range("E16").Value = "  +1.43%  "
$ws.Range("D17").Value = "66.823.74"
$ws.Range("E17").Value = "  +0.94%  "
$ws.Range("E18").Value = "  +1.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.60"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.39"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +3.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.20"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "429.64"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.623"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.98"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000123"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +4.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.40"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +5.78%  "
$ws.Range("E28").Value = "  +4.96%  "
$ws.Range("E29").Value = "  +0.94%  "
$ws.Range("E30").Value = "  -0.15%  "
$ws.Range("D31").Value = "3.632.44"
$ws.Range("E31").Value = "  +2.10%  "
$ws.Range("E32").Value = "  +1.53%  "
$ws.Range("E33").Value = "  +3.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.53"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.92"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.91%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.71"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.79%  "
$ws.Range("E38").Value = "  -1.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "176.98"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0865"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.04%  "
$ws.Range("E41").Value = "  +1.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.903"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.79%  "
$ws.Range("E43").Value = "  -0.84%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.58"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +8.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.999"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "25.22"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.19"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "24.03"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.23"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.962"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.02%  "
